# Reschedule the "next team meeting" date: Tuesday 18 March -> Tuesday 19 March,
# and relocate Word's "_GoBack" (last-edit-position) bookmark to sit right after the
# newly typed "9", matching where the author's cursor was when they made the edit.

$d = $word.ActiveDocument

# Locate the unique phrase containing the date we need to change.
$find = $d.Content
$found = $find.Find.Execute("Tuesday 18 March", $true, $false, $false, $false, $false, `
                             $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find 'Tuesday 18 March' in the document"
}

$matchStart = $find.Start

# "Tuesday 18 March"
#  0123456789...
# The "8" in "18" sits 9 characters after the start of the match.
$digitRange = $d.Range($matchStart + 9, $matchStart + 10)
$digitRange.Text = "9"

# Word re-anchors the hidden "_GoBack" bookmark to the location of the most recent edit;
# here that is immediately after the freshly-typed "9".
$goBackRange = $d.Range($matchStart + 10, $matchStart + 10)
$d.Bookmarks.Add("_GoBack", $goBackRange)
